# Apply cryptos list update (prices, volumes, and a coin-order swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.751.20"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.601.32"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.52"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.81"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.623.72"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +6.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.368"
$ws.Range("E13").Value = "  +7.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.062.68"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.56"
$ws.Range("E15").Value = "  +6.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.726.56"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.607.40"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.62"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.75"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.64"
$ws.Range("E21").Value = "  +4.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("E22").Value = "  +10.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.522"
$ws.Range("E24").Value = "  +15.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.43"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.59"
$ws.Range("E28").Value = "  +4.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.87"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.39"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.13"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.18"
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.918"
$ws.Range("E37").Value = "  +4.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.75"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.72"
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "292.25"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.12"
$ws.Range("E43").Value = "  +12.31%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0540"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.66"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.81"
$ws.Range("E50").Value = "  +7.49%  "
$ws.Range("E51").Value = "  +2.81%  "
